# Updated symbol list on Wed Jan 11 23:14:23 UTC 2023 with GitHub Actions
# Refreshes Price/Volume(1h)/Hora columns for each coin row, and swaps the
# KickToken/BKEXToken rows (41/42) to match the new ranking order.
# Values are written with a leading "'" (quote-prefix) so Excel stores them
# as literal text (matching the source inlineStr cells) instead of
# reinterpreting numeric-looking strings ("280.54", "23") or percentages
# ("0.84%") as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'280.54"
$ws.Range("E2").Value = "'0.84%"
$ws.Range("G2").Value = "'23"

$ws.Range("D3").Value = "'27.52"
$ws.Range("E3").Value = "'1.24%"
$ws.Range("G3").Value = "'23"

$ws.Range("D4").Value = "'4.846"
$ws.Range("E4").Value = "'-0.76%"
$ws.Range("G4").Value = "'23"

$ws.Range("D5").Value = "'0.06430"
$ws.Range("E5").Value = "'0.26%"
$ws.Range("G5").Value = "'23"

$ws.Range("D6").Value = "'7.081"
$ws.Range("E6").Value = "'1.81%"
$ws.Range("G6").Value = "'23"

$ws.Range("D7").Value = "'1.278"
$ws.Range("E7").Value = "'2.31%"
$ws.Range("G7").Value = "'23"

$ws.Range("D8").Value = "'0.9011"
$ws.Range("E8").Value = "'2.22%"
$ws.Range("G8").Value = "'23"

$ws.Range("D9").Value = "'0.1545"
$ws.Range("E9").Value = "'1.23%"
$ws.Range("G9").Value = "'23"

$ws.Range("D10").Value = "'0.06520"
$ws.Range("E10").Value = "'29.92%"
$ws.Range("G10").Value = "'23"

$ws.Range("D11").Value = "'0.07490"
$ws.Range("E11").Value = "'-0.50%"
$ws.Range("G11").Value = "'23"

$ws.Range("D12").Value = "'0.02929"
$ws.Range("E12").Value = "'3.23%"
$ws.Range("G12").Value = "'23"

$ws.Range("D13").Value = "'0.08991"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("G13").Value = "'23"

$ws.Range("D14").Value = "'0.001589"
$ws.Range("E14").Value = "'1.65%"
$ws.Range("G14").Value = "'23"

$ws.Range("D15").Value = "'0.0006430"
$ws.Range("E15").Value = "'0.46%"
$ws.Range("G15").Value = "'23"

$ws.Range("D16").Value = "'0.006054"
$ws.Range("E16").Value = "'1.06%"
$ws.Range("G16").Value = "'23"

$ws.Range("D17").Value = "'3.488"
$ws.Range("E17").Value = "'0.84%"
$ws.Range("G17").Value = "'23"

$ws.Range("D18").Value = "'3.303"
$ws.Range("E18").Value = "'-0.36%"
$ws.Range("G18").Value = "'23"

$ws.Range("D19").Value = "'2.224"
$ws.Range("E19").Value = "'-2.11%"
$ws.Range("G19").Value = "'23"

$ws.Range("E20").Value = "'-1.29%"
$ws.Range("G20").Value = "'23"

$ws.Range("D21").Value = "'0.1353"
$ws.Range("E21").Value = "'1.85%"
$ws.Range("G21").Value = "'23"

$ws.Range("D22").Value = "'3.895"
$ws.Range("E22").Value = "'-0.27%"
$ws.Range("G22").Value = "'23"

$ws.Range("D23").Value = "'0.04400"
$ws.Range("E23").Value = "'-0.67%"
$ws.Range("G23").Value = "'23"

$ws.Range("D24").Value = "'0.1502"
$ws.Range("E24").Value = "'8.80%"
$ws.Range("G24").Value = "'23"

$ws.Range("D25").Value = "'0.001174"
$ws.Range("E25").Value = "'-0.24%"
$ws.Range("G25").Value = "'23"

$ws.Range("D26").Value = "'0.004300"
$ws.Range("E26").Value = "'11.37%"
$ws.Range("G26").Value = "'23"

$ws.Range("G27").Value = "'23"

$ws.Range("D28").Value = "'0.0001178"
$ws.Range("E28").Value = "'-1.86%"
$ws.Range("G28").Value = "'23"

$ws.Range("D29").Value = "'0.0001656"
$ws.Range("G29").Value = "'23"

$ws.Range("G30").Value = "'23"

$ws.Range("G31").Value = "'23"

$ws.Range("G32").Value = "'23"

$ws.Range("G33").Value = "'23"

$ws.Range("G34").Value = "'23"

$ws.Range("G35").Value = "'23"

$ws.Range("G36").Value = "'23"

$ws.Range("G37").Value = "'23"

$ws.Range("G38").Value = "'23"

$ws.Range("G39").Value = "'23"

$ws.Range("D40").Value = "'0.04089"
$ws.Range("E40").Value = "'-1.28%"
$ws.Range("G40").Value = "'23"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1412"
$ws.Range("E41").Value = "'19.97%"
$ws.Range("G41").Value = "'23"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006627"
$ws.Range("E42").Value = "'-2.46%"
$ws.Range("G42").Value = "'23"

$ws.Range("D43").Value = "'0.002076"
$ws.Range("E43").Value = "'-13.21%"
$ws.Range("G43").Value = "'23"

$ws.Range("D44").Value = "'0.01103"
$ws.Range("E44").Value = "'-1.97%"
$ws.Range("G44").Value = "'23"

$ws.Range("D45").Value = "'0.00005549"
$ws.Range("E45").Value = "'6.74%"
$ws.Range("G45").Value = "'23"

$ws.Range("D46").Value = "'1.561"
$ws.Range("E46").Value = "'5.28%"
$ws.Range("G46").Value = "'23"

$ws.Range("E47").Value = "'-8.78%"
$ws.Range("G47").Value = "'23"

$ws.Range("G48").Value = "'23"

$ws.Range("G49").Value = "'23"

$ws.Range("G50").Value = "'23"

$ws.Range("G51").Value = "'23"
